$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for new rows 252-255 (dates 44326-44329), mirroring the pattern of row 251
$rows = @(
    @{ Row = 252; A = 44326; B = 0; C = 1; D = 28.87669650591972 },
    @{ Row = 253; A = 44327; B = 0; C = 1; D = 28.87669650591972 },
    @{ Row = 254; A = 44328; B = 0; C = 1; D = 28.87669650591972 },
    @{ Row = 255; A = 44329; B = 0; C = 1; D = 28.87669650591972 }
)

# Copy the style of A251 (date format) down to the new date cells
$styleSource = $ws.Range("A251")

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D

    $styleSource.Copy()
    $ws.Cells.Item($r.Row, 1).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false
